$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 29.41996
$ws.Range("H2").Value = 88.25988000000001
$ws.Range("I2").Value = 0.6207199949605289
$ws.Range("J2").Value = 0.66829493802317
$ws.Range("M2").Value = 293.7577056666667
$ws.Range("N2").Value = 881.273117
$ws.Range("O2").Value = 0.9369756110667984
$ws.Range("P2").Value = 0.940660486426629
$ws.Range("Q2").Value = 8642.339950405107
$ws.Range("R2").Value = 77781.05955364596
$ws.Range("S2").Value = 0.5815994965795216
$ws.Range("T2").Value = 0.628638641477329
# Row 3
$ws.Range("G3").Value = 29.41996
$ws.Range("H3").Value = 88.25988000000001
$ws.Range("I3").Value = 0.6207199949605289
$ws.Range("J3").Value = 0.66829493802317
$ws.Range("O3").Value = 0.01525159481997056
$ws.Range("P3").Value = 0.01531157527761154
$ws.Range("Q3").Value = 140.67545159468
$ws.Range("R3").Value = 1266.07906435212
$ws.Range("S3").Value = 0.009466969859792153
$ws.Range("T3").Value = 0.01023264825118851
# Row 4
$ws.Range("G4").Value = 29.41996
$ws.Range("H4").Value = 88.25988000000001
$ws.Range("I4").Value = 0.6207199949605289
$ws.Range("J4").Value = 0.66829493802317
$ws.Range("M4").Value = 7.906212666666666
$ws.Range("N4").Value = 23.718638
$ws.Range("O4").Value = 0.02521781829607561
$ws.Range("P4").Value = 0.02531699325449539
$ws.Range("Q4").Value = 232.6004604048267
$ws.Range("R4").Value = 2093.40414364344
$ws.Range("S4").Value = 0.01565320404565559
$ws.Range("T4").Value = 0.01691921843794601
# Row 5
$ws.Range("G5").Value = 29.41996
$ws.Range("H5").Value = 88.25988000000001
$ws.Range("I5").Value = 0.6207199949605289
$ws.Range("J5").Value = 0.66829493802317
$ws.Range("M5").Value = 3.6844455
$ws.Range("N5").Value = 7.368891000000001
$ws.Range("O5").Value = 0.01175198303639443
$ws.Range("P5").Value = 0.007865466969060864
$ws.Range("Q5").Value = 108.39623923218
$ws.Range("R5").Value = 650.3774353930801
$ws.Range("S5").Value = 0.007294690851126975
$ws.Range("T5").Value = 0.005256451760611821
# Row 6
$ws.Range("G6").Value = 29.41996
$ws.Range("H6").Value = 88.25988000000001
$ws.Range("I6").Value = 0.6207199949605289
$ws.Range("J6").Value = 0.66829493802317
$ws.Range("M6").Value = 3.386921
$ws.Range("N6").Value = 10.160763
$ws.Range("O6").Value = 0.01080299278076119
$ws.Range("P6").Value = 0.01084547807220323
$ws.Range("Q6").Value = 99.64308034315999
$ws.Range("R6").Value = 896.7877230884401
$ws.Range("S6").Value = 0.006705633624432718
$ws.Range("T6").Value = 0.007247978096094709
# Row 7
$ws.Range("I7").Value = 0.1515698101047853
$ws.Range("J7").Value = 0.1631868437822795
$ws.Range("M7").Value = 293.7577056666667
$ws.Range("N7").Value = 881.273117
$ws.Range("O7").Value = 0.9369756110667984
$ws.Range("P7").Value = 0.940660486426629
$ws.Range("Q7").Value = 2110.320008665418
$ws.Range("R7").Value = 18992.88007798876
$ws.Range("S7").Value = 0.1420172154422098
$ws.Range("T7").Value = 0.1535034158506653
# Row 8
$ws.Range("I8").Value = 0.1515698101047853
$ws.Range("J8").Value = 0.1631868437822795
$ws.Range("O8").Value = 0.01525159481997056
$ws.Range("P8").Value = 0.01531157527761154
$ws.Range("R8").Value = 309.156084739461
$ws.Range("S8").Value = 0.002311681330658064
$ws.Range("T8").Value = 0.002498647642888208
# Row 9
$ws.Range("I9").Value = 0.1515698101047853
$ws.Range("J9").Value = 0.1631868437822795
$ws.Range("M9").Value = 7.906212666666666
$ws.Range("N9").Value = 23.718638
$ws.Range("O9").Value = 0.02521781829607561
$ws.Range("P9").Value = 0.02531699325449539
$ws.Range("Q9").Value = 56.79728041640911
$ws.Range("R9").Value = 511.175523747682
$ws.Range("S9").Value = 0.00382225993039316
$ws.Range("T9").Value = 0.004131400223258363
# Row 10
$ws.Range("I10").Value = 0.1515698101047853
$ws.Range("J10").Value = 0.1631868437822795
$ws.Range("M10").Value = 3.6844455
$ws.Range("N10").Value = 7.368891000000001
$ws.Range("O10").Value = 0.01175198303639443
$ws.Range("P10").Value = 0.007865466969060864
$ws.Range("Q10").Value = 26.4686131103915
$ws.Range("R10").Value = 158.811678662349
$ws.Range("S10").Value = 0.001781245837180962
$ws.Range("T10").Value = 0.001283540729554815
# Row 11
$ws.Range("I11").Value = 0.1515698101047853
$ws.Range("J11").Value = 0.1631868437822795
$ws.Range("M11").Value = 3.386921
$ws.Range("N11").Value = 10.160763
$ws.Range("O11").Value = 0.01080299278076119
$ws.Range("P11").Value = 0.01084547807220323
$ws.Range("Q11").Value = 24.33123290450633
$ws.Range("R11").Value = 218.981096140557
$ws.Range("S11").Value = 0.00163740756434334
$ws.Range("T11").Value = 0.001769839335912767
# Row 12
$ws.Range("G12").Value = 0.3873096666666667
$ws.Range("H12").Value = 1.161929
$ws.Range("I12").Value = 0.008171692087327698
$ws.Range("J12").Value = 0.008798009571759262
$ws.Range("M12").Value = 293.7577056666667
$ws.Range("N12").Value = 881.273117
$ws.Range("O12").Value = 0.9369756110667984
$ws.Range("P12").Value = 0.940660486426629
$ws.Range("Q12").Value = 113.7751990625214
$ws.Range("R12").Value = 1023.976791562693
$ws.Range("S12").Value = 0.007656676186973591
$ws.Range("T12").Value = 0.008275939963357206
# Row 13
$ws.Range("G13").Value = 0.3873096666666667
$ws.Range("H13").Value = 1.161929
$ws.Range("I13").Value = 0.008171692087327698
$ws.Range("J13").Value = 0.008798009571759262
$ws.Range("O13").Value = 0.01525159481997056
$ws.Range("P13").Value = 0.01531157527761154
$ws.Range("Q13").Value = 1.851972683352333
$ws.Range("R13").Value = 16.667754150171
$ws.Range("S13").Value = 0.0001246313367094815
$ws.Range("T13").Value = 0.0001347113858511389
# Row 14
$ws.Range("G14").Value = 0.3873096666666667
$ws.Range("H14").Value = 1.161929
$ws.Range("I14").Value = 0.008171692087327698
$ws.Range("J14").Value = 0.008798009571759262
$ws.Range("M14").Value = 7.906212666666666
$ws.Range("N14").Value = 23.718638
$ws.Range("O14").Value = 0.02521781829607561
$ws.Range("P14").Value = 0.02531699325449539
$ws.Range("Q14").Value = 3.062152592522444
$ws.Range("R14").Value = 27.559373332702
$ws.Range("S14").Value = 0.0002060722462297087
$ws.Range("T14").Value = 0.0002227391489812151
# Row 15
$ws.Range("G15").Value = 0.3873096666666667
$ws.Range("H15").Value = 1.161929
$ws.Range("I15").Value = 0.008171692087327698
$ws.Range("J15").Value = 0.008798009571759262
$ws.Range("M15").Value = 3.6844455
$ws.Range("N15").Value = 7.368891000000001
$ws.Range("O15").Value = 0.01175198303639443
$ws.Range("P15").Value = 0.007865466969060864
$ws.Range("Q15").Value = 1.4270213584565
$ws.Range("R15").Value = 8.562128150739001
$ws.Range("S15").Value = 0.00009603358678891374
$ws.Range("T15").Value = 0.0000692004536801538
# Row 16
$ws.Range("G16").Value = 0.3873096666666667
$ws.Range("H16").Value = 1.161929
$ws.Range("I16").Value = 0.008171692087327698
$ws.Range("J16").Value = 0.008798009571759262
$ws.Range("M16").Value = 3.386921
$ws.Range("N16").Value = 10.160763
$ws.Range("O16").Value = 0.01080299278076119
$ws.Range("P16").Value = 0.01084547807220323
$ws.Range("Q16").Value = 1.311787243536333
$ws.Range("R16").Value = 11.806085191827
$ws.Range("S16").Value = 0.00008827873062600451
$ws.Range("T16").Value = 0.00009541861988954924
# Row 17
$ws.Range("G17").Value = 10.122265
$ws.Range("H17").Value = 20.24453
$ws.Range("I17").Value = 0.2135656295858028
$ws.Range("J17").Value = 0.153289545846405
$ws.Range("M17").Value = 293.7577056666667
$ws.Range("N17").Value = 881.273117
$ws.Range("O17").Value = 0.9369756110667984
$ws.Range("P17").Value = 0.940660486426629
$ws.Range("Q17").Value = 2973.493342550001
$ws.Range("R17").Value = 17840.96005530001
$ws.Range("S17").Value = 0.2001057862840231
$ws.Range("T17").Value = 0.1441934187599964
# Row 18
$ws.Range("G18").Value = 10.122265
$ws.Range("H18").Value = 20.24453
$ws.Range("I18").Value = 0.2135656295858028
$ws.Range("J18").Value = 0.153289545846405
$ws.Range("O18").Value = 0.01525159481997056
$ws.Range("P18").Value = 0.01531157527761154
$ws.Range("Q18").Value = 48.400956358745
$ws.Range("R18").Value = 290.40573815247
$ws.Range("S18").Value = 0.003257216449914581
$ws.Range("T18").Value = 0.002347104420498116
# Row 19
$ws.Range("G19").Value = 10.122265
$ws.Range("H19").Value = 20.24453
$ws.Range("I19").Value = 0.2135656295858028
$ws.Range("J19").Value = 0.153289545846405
$ws.Range("M19").Value = 7.906212666666666
$ws.Range("N19").Value = 23.718638
$ws.Range("O19").Value = 0.02521781829607561
$ws.Range("P19").Value = 0.02531699325449539
$ws.Range("Q19").Value = 80.02877975835665
$ws.Range("R19").Value = 480.1726785501399
$ws.Range("S19").Value = 0.005385659241181766
$ws.Range("T19").Value = 0.003880830398178097
# Row 20
$ws.Range("G20").Value = 10.122265
$ws.Range("H20").Value = 20.24453
$ws.Range("I20").Value = 0.2135656295858028
$ws.Range("J20").Value = 0.153289545846405
$ws.Range("M20").Value = 3.6844455
$ws.Range("N20").Value = 7.368891000000001
$ws.Range("O20").Value = 0.01175198303639443
$ws.Range("P20").Value = 0.007865466969060864
$ws.Range("Q20").Value = 37.2949337290575
$ws.Range("R20").Value = 149.17973491623
$ws.Range("S20").Value = 0.002509819656049252
$ws.Range("T20").Value = 0.001205693859557239
# Row 21
$ws.Range("G21").Value = 10.122265
$ws.Range("H21").Value = 20.24453
$ws.Range("I21").Value = 0.2135656295858028
$ws.Range("J21").Value = 0.153289545846405
$ws.Range("M21").Value = 3.386921
$ws.Range("N21").Value = 10.160763
$ws.Range("O21").Value = 0.01080299278076119
$ws.Range("P21").Value = 0.01084547807220323
$ws.Range("Q21").Value = 34.28331189606499
$ws.Range("R21").Value = 205.69987137639
$ws.Range("S21").Value = 0.002307147954634147
$ws.Range("T21").Value = 0.001662498408175178
# Row 22
$ws.Range("G22").Value = 0.2830933333333334
$ws.Range("H22").Value = 0.84928
$ws.Range("I22").Value = 0.005972873261555284
$ws.Range("J22").Value = 0.006430662776386256
$ws.Range("M22").Value = 293.7577056666667
$ws.Range("N22").Value = 881.273117
$ws.Range("O22").Value = 0.9369756110667984
$ws.Range("P22").Value = 0.940660486426629
$ws.Range("Q22").Value = 83.1608480895289
$ws.Range("R22").Value = 748.44763280576
$ws.Range("S22").Value = 0.005596436574070303
$ws.Range("T22").Value = 0.006049070375281112
# Row 23
$ws.Range("G23").Value = 0.2830933333333334
$ws.Range("H23").Value = 0.84928
$ws.Range("I23").Value = 0.005972873261555284
$ws.Range("J23").Value = 0.006430662776386256
$ws.Range("O23").Value = 0.01525159481997056
$ws.Range("P23").Value = 0.01531157527761154
$ws.Range("Q23").Value = 1.353648424746667
$ws.Range("R23").Value = 12.18283582272
$ws.Range("S23").Value = 0.00009109584289627721
$ws.Range("T23").Value = 0.00009846357718557261
# Row 24
$ws.Range("G24").Value = 0.2830933333333334
$ws.Range("H24").Value = 0.84928
$ws.Range("I24").Value = 0.005972873261555284
$ws.Range("J24").Value = 0.006430662776386256
$ws.Range("M24").Value = 7.906212666666666
$ws.Range("N24").Value = 23.718638
$ws.Range("O24").Value = 0.02521781829607561
$ws.Range("P24").Value = 0.02531699325449539
$ws.Range("Q24").Value = 2.238196097848889
$ws.Range("R24").Value = 20.14376488064
$ws.Range("S24").Value = 0.0001506228326153896
$ws.Range("T24").Value = 0.0001628050461317054
# Row 25
$ws.Range("G25").Value = 0.2830933333333334
$ws.Range("H25").Value = 0.84928
$ws.Range("I25").Value = 0.005972873261555284
$ws.Range("J25").Value = 0.006430662776386256
$ws.Range("M25").Value = 3.6844455
$ws.Range("N25").Value = 7.368891000000001
$ws.Range("O25").Value = 0.01175198303639443
$ws.Range("P25").Value = 0.007865466969060864
$ws.Range("Q25").Value = 1.04304195808
$ws.Range("R25").Value = 6.258251748480001
$ws.Range("S25").Value = 0.0000701931052483316
$ws.Range("T25").Value = 0.00005058016565683533
# Row 26
$ws.Range("G26").Value = 0.2830933333333334
$ws.Range("H26").Value = 0.84928
$ws.Range("I26").Value = 0.005972873261555284
$ws.Range("J26").Value = 0.006430662776386256
$ws.Range("M26").Value = 3.386921
$ws.Range("N26").Value = 10.160763
$ws.Range("O26").Value = 0.01080299278076119
$ws.Range("P26").Value = 0.01084547807220323
$ws.Range("Q26").Value = 0.9588147556266666
$ws.Range("R26").Value = 8.62933280064
$ws.Range("S26").Value = 0.0000645249067249833
$ws.Range("T26").Value = 0.00006974361213103071
